$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.573.76'
$ws.Range("E2").Value = '  -1.36%  '
$ws.Range("D3").Value = '1.846.17'
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -1.16%  '
$ws.Range("D5").Value = '''333.38'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").Value = '''1.004'
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").Value = '''0.4638'
$ws.Range("E7").Value = '  -1.31%  '
$ws.Range("D8").Value = '''0.3856'
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("D9").Value = '''46.25'
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").Value = '''0.07892'
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("D11").Value = '''0.9943'
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("D12").Value = '''21.45'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '1.851.26'
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").Value = '''5.916'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").Value = '''7.104'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '''1.006'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '''88.93'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").Value = '''0.06644'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").Value = '''17.05'
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("E21").Value = '  -1.06%  '
$ws.Range("D22").Value = '27.577.86'
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").Value = '''5.377'
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").Value = '''10.90'
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").Value = '''2.304'
$ws.Range("E25").Value = '  -2.45%  '
$ws.Range("D26").Value = '''158.16'
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").Value = '''19.51'
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").Value = '''2.097'
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("D29").Value = '''5.400'
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("E30").Value = '  -1.34%  '
$ws.Range("D31").Value = '''0.9760'
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").Value = '''0.09403'
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").Value = '''3.585'
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("D34").Value = '''5.283'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").Value = '''1.340'
$ws.Range("E35").Value = '  -1.59%  '
$ws.Range("D36").Value = '''0.06022'
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").Value = '''0.02225'
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D38").Value = '''8.282'
$ws.Range("E38").Value = '  +0.98%  '
$ws.Range("D39").Value = '''1.182'
$ws.Range("E39").Value = '  -2.87%  '
$ws.Range("D40").Value = '''0.5887'
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("D41").Value = '''0.1862'
$ws.Range("E41").Value = '  -2.06%  '
$ws.Range("D42").Value = '''10.28'
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("D43").Value = '''1.253'
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").Value = '''0.5574'
$ws.Range("E44").Value = '  -2.03%  '
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = '''1.897'
$ws.Range("E46").Value = '  -2.33%  '
$ws.Range("D47").Value = '''0.06680'
$ws.Range("E47").Value = '  -2.82%  '
$ws.Range("D48").Value = '''110.71'
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("D49").Value = '''1.052'
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("D50").Value = '''1.003'
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("E51").Value = '  -1.62%  '
